$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33; this shifts the existing rows 33..138 down to 34..139
$ws.Rows.Item(33).Insert()

# Populate the new row 33 with a fresh weekly price entry
# (same market / category / variety / quality as the former row 33, new date & prices)
$ws.Cells.Item(33, 1).Value = 4
$ws.Cells.Item(33, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(33, 3).Value = "Los Lagos"
$ws.Cells.Item(33, 4).Value = 44519
$ws.Cells.Item(33, 5).Value = 10
$ws.Cells.Item(33, 6).Value = 100112028
$ws.Cells.Item(33, 7).Value = "Sandia"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 700
$ws.Cells.Item(33, 11).Value = 1200
$ws.Cells.Item(33, 12).Value = 1200
$ws.Cells.Item(33, 13).Value = 1200
$ws.Cells.Item(33, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(33, 15).Value = "Perú"
$ws.Cells.Item(33, 16).Value = 1200
$ws.Cells.Item(33, 17).Value = 1
$ws.Cells.Item(33, 18).Value = "Hortaliza"
